$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 265, shifting existing rows 265-276 down to 266-277.
$ws.Rows.Item(265).Insert()

# Populate the new row 265 with the new record's data.
$ws.Range("A265").Value = 6
$ws.Range("B265").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C265").Value = 'Metropolitana'
$ws.Range("D265").Value = 44939
$ws.Range("E265").Value = 13
$ws.Range("F265").Value = 100112001
$ws.Range("G265").Value = 'Berenjena'
$ws.Range("H265").Value = 'Sin especificar'
$ws.Range("I265").Value = 'Primera'
$ws.Range("J265").Value = 380
$ws.Range("K265").Value = 9000
$ws.Range("L265").Value = 9000
$ws.Range("M265").Value = 9000
$ws.Range("N265").Value = '$/caja 40 unidades'
$ws.Range("O265").Value = 'Región Metropolitana'
$ws.Range("P265").Value = 225
$ws.Range("Q265").Value = 40
$ws.Range("R265").Value = 'Hortaliza'
